$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix erroneous IFRS consolidated figures (values were off by a scale/mapping error) ---

# Row 2: correct financial figures
$ws.Range("D2").Value = 539
$ws.Range("E2").Value = -21
$ws.Range("F2").Value = -21
$ws.Range("G2").Value = -23
$ws.Range("H2").Value = -20
$ws.Range("I2").Value = -20
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 516
$ws.Range("L2").Value = 244
$ws.Range("M2").Value = 272
$ws.Range("N2").Value = 270
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 71
$ws.Range("Q2").Value = -5
$ws.Range("R2").Value = -29
$ws.Range("S2").Value = 28
$ws.Range("T2").Value = 25
$ws.Range("U2").Value = -30
$ws.Range("V2").Value = 168
$ws.Range("W2").Value = -3.83
$ws.Range("X2").Value = -3.66
$ws.Range("Y2").Value = -6.88
$ws.Range("Z2").Value = -3.9
$ws.Range("AA2").Value = 89.69
$ws.Range("AB2").Value = 290.15
$ws.Range("AC2").Value = -138
$ws.Range("AD2").Value = -9.6
$ws.Range("AE2").Value = 2077
$ws.Range("AF2").Value = 0.64
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = -0.15
$ws.Range("AJ2").Value = 13621761

# Row 3: correct financial figures
$ws.Range("D3").Value = 482
$ws.Range("E3").Value = -13
$ws.Range("F3").Value = -13
$ws.Range("G3").Value = -17
$ws.Range("H3").Value = -17
$ws.Range("I3").Value = -17
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 453
$ws.Range("L3").Value = 167
$ws.Range("M3").Value = 286
$ws.Range("N3").Value = 283
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 82
$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = -27
$ws.Range("S3").Value = -44
$ws.Range("T3").Value = 27
$ws.Range("U3").Value = -19
$ws.Range("V3").Value = 97
$ws.Range("W3").Value = -2.79
$ws.Range("X3").Value = -3.52
$ws.Range("Y3").Value = -6.06
$ws.Range("Z3").Value = -3.5
$ws.Range("AA3").Value = 58.52
$ws.Range("AB3").Value = 251.38
$ws.Range("AC3").Value = -108
$ws.Range("AD3").Value = -26.3
$ws.Range("AE3").Value = 1859
$ws.Range("AF3").Value = 1.53
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = -0.18
$ws.Range("AJ3").Value = 15894485

# Row 4: correct financial figures
$ws.Range("D4").Value = 454
$ws.Range("E4").Value = -7
$ws.Range("F4").Value = -7
$ws.Range("G4").Value = -11
$ws.Range("H4").Value = -11
$ws.Range("I4").Value = -11
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 441
$ws.Range("L4").Value = 161
$ws.Range("M4").Value = 280
$ws.Range("N4").Value = 278
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 86
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = -19
$ws.Range("T4").Value = 8
$ws.Range("U4").Value = 18
$ws.Range("V4").Value = 69
$ws.Range("W4").Value = -1.56
$ws.Range("X4").Value = -2.35
$ws.Range("Y4").Value = -3.74
$ws.Range("Z4").Value = -2.39
$ws.Range("AA4").Value = 57.38
$ws.Range("AB4").Value = 236
$ws.Range("AC4").Value = -62
$ws.Range("AD4").Value = -165.75
$ws.Range("AE4").Value = 1750
$ws.Range("AF4").Value = 5.83
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = -0.29
$ws.Range("AJ4").Value = 16554309

# Row 5: correct financial figures
$ws.Range("D5").Value = 461
$ws.Range("E5").Value = -5
$ws.Range("F5").Value = -5
$ws.Range("G5").Value = -12
$ws.Range("H5").Value = -27
$ws.Range("I5").Value = -27
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 430
$ws.Range("L5").Value = 137
$ws.Range("M5").Value = 293
$ws.Range("N5").Value = 291
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 100
$ws.Range("Q5").Value = 24
$ws.Range("R5").Value = -27
$ws.Range("S5").Value = 22
$ws.Range("T5").Value = 23
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 47
$ws.Range("W5").Value = -0.99
$ws.Range("X5").Value = -5.79
$ws.Range("Y5").Value = -9.390000000000001
$ws.Range("Z5").Value = -6.12
$ws.Range("AA5").Value = 46.72
$ws.Range("AB5").Value = 203.84
$ws.Range("AC5").Value = -154
$ws.Range("AD5").Value = -14.34
$ws.Range("AE5").Value = 1545
$ws.Range("AF5").Value = 1.43
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = -0.11
$ws.Range("AJ5").Value = 19486857

# Row 6: correct financial figures
$ws.Range("D6").Value = 413
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 2
$ws.Range("K6").Value = 433
$ws.Range("L6").Value = 139
$ws.Range("M6").Value = 294
$ws.Range("N6").Value = 292
$ws.Range("P6").Value = 100
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = -68
$ws.Range("S6").Value = 24
$ws.Range("T6").Value = 38
$ws.Range("U6").Value = -28
$ws.Range("V6").Value = 72
$ws.Range("W6").Value = 2.03
$ws.Range("X6").Value = 0.47
$ws.Range("Y6").Value = 0.7
$ws.Range("Z6").Value = 0.45
$ws.Range("AA6").Value = 47.34
$ws.Range("AB6").Value = 215.83
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 259.69
$ws.Range("AE6").Value = 1549
$ws.Range("AF6").Value = 1.71
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 1.47
$ws.Range("AJ6").Value = 19486857

# Row 7: remove stale/incorrect figures (keep only A/B/C identifying columns)
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: remove stale/incorrect figures (keep only A/B/C identifying columns)
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: remove stale/incorrect figures (keep only A/B/C identifying columns)
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
